$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.838884592056274
$ws.Range("B1").Value = 4.521198749542236
$ws.Range("C1").Value = 3.215765476226807
$ws.Range("D1").Value = 2.477200508117676
$ws.Range("E1").Value = 2.292156219482422
